# Add an "active" column to the front of the Debts and Fixed Assets sheets.
$wb = $excel.ActiveWorkbook

$wsDebts = $wb.Worksheets.Item("Debts")
$wsFixed = $wb.Worksheets.Item("Fixed Assets")

# Debts: insert a new column A and label it "active"; existing columns shift right.
$wsDebts.Range("A1").EntireColumn.Insert()
$wsDebts.Range("A1").Value = "active"

# Fixed Assets: insert a new column A and label it "active"; existing columns shift right.
$wsFixed.Range("A1").EntireColumn.Insert()
$wsFixed.Range("A1").Value = "active"

# Restore/update each sheet's selection.
[void]$wsFixed.Range("C9").Select()

# Make Debts the active sheet/tab, with its own selection.
[void]$wsDebts.Activate()
[void]$wsDebts.Range("B10").Select()
